$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = "/etakɾã/ [e.ta.'kɾã] (uno; solo) (test comment) (test comment 2){4}; /etakrã/"
$ws.Range("H3").Value = "<peteĩ>(uno){Guasch1962:670} (Test comment 3)"

$ws.Range("G9").Select()
